$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7241352796554565
$ws.Range("B1").Value = 1.186225533485413
$ws.Range("C1").Value = 2.383294820785522
$ws.Range("D1").Value = 3.465487241744995
$ws.Range("E1").Value = 3.12832236289978
